$d = $word.ActiveDocument

# --- Step 1: fix the wording -----------------------------------------
# "A CUENTA POR <TIPO_ABONO>" -> "A CUENTA <TIPO_ABONO>"
# (drops "POR ", keeps a single space before the placeholder). The whole
# phrase currently lives in one run, so the replace keeps that run's
# formatting (bold, purple, sz 20) intact.
$findRng = $d.Content
$ok = $findRng.Find.Execute("A CUENTA POR <TIPO_ABONO>", $true, $false, $false, $false, $false, `
                             $true, 1, $false, "A CUENTA <TIPO_ABONO>", 2)
if (-not $ok) {
    throw "Could not find the target phrase 'A CUENTA POR <TIPO_ABONO>'"
}

# --- Step 2: locate the freshly-written text and split it into 3 runs -
# Target markup needs three separate <w:r> runs (identical rPr):
#   "A CUENTA"  |  " "  |  "<TIPO_ABONO>"
$locate = $d.Content
$locate.Find.Execute("A CUENTA <TIPO_ABONO>")
$start = $locate.Start

# Offsets within "A CUENTA <TIPO_ABONO>" where the run boundaries fall.
$splitAfterWord = $start + 8   # right after "A CUENTA"
$splitAfterSpace = $start + 9  # right after the following space

# Toggling a character property off then back on forces Word to break
# the run at that boundary without altering the final formatting, since
# the net effect on the run properties is a no-op.
$part1 = $d.Range($start, $splitAfterWord)
$part1.Font.Bold = $false
$part1.Font.Bold = $true

$part2 = $d.Range($splitAfterWord, $splitAfterSpace)
$part2.Font.Bold = $false
$part2.Font.Bold = $true

Write-Output "Split 'A CUENTA <TIPO_ABONO>' into three runs successfully"
